$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column price/volume cells are text in the source data (e.g. "37.708.14",
# "1.10"), not real numbers -- force Text format first so Excel does not
# "smart" convert/truncate them when the value looks numeric, then restore
# General formatting so the cell style matches the original.
$textCells = @('D2', 'D3', 'D5', 'D6', 'D7', 'D10', 'D12', 'D13', 'D14', 'D15', 'D16', 'D17', 'D18', 'D19', 'D20', 'D21', 'D22', 'D24', 'D25', 'D26', 'D27', 'D28', 'D29', 'D30', 'D31', 'D33', 'D34', 'D36', 'D37', 'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D50', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '37.655.25'
$ws.Range('E2').Value = '  -0.20%  '
$ws.Range('D3').Value = '2.035.30'
$ws.Range('E3').Value = '  +0.56%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '227.12'
$ws.Range('E5').Value = '  +0.08%  '
$ws.Range('D6').Value = '0.606'
$ws.Range('E6').Value = '  -1.00%  '
$ws.Range('D7').Value = '59.38'
$ws.Range('E7').Value = '  -0.34%  '
$ws.Range('E9').Value = '  -2.33%  '
$ws.Range('D10').Value = '0.0834'
$ws.Range('E10').Value = '  +2.47%  '
$ws.Range('E11').Value = '  -0.19%  '
$ws.Range('D12').Value = '2.336.49'
$ws.Range('E12').Value = '  +0.57%  '
$ws.Range('D13').Value = '14.42'
$ws.Range('E13').Value = '  -0.88%  '
$ws.Range('D14').Value = '20.98'
$ws.Range('E14').Value = '  +0.18%  '
$ws.Range('D15').Value = '5.43'
$ws.Range('E15').Value = '  +4.17%  '
$ws.Range('D16').Value = '0.769'
$ws.Range('E16').Value = '  +2.47%  '
$ws.Range('D17').Value = '2.045.93'
$ws.Range('E17').Value = '  -0.13%  '
$ws.Range('D18').Value = '37.629.18'
$ws.Range('E18').Value = '  -0.23%  '
$ws.Range('B19').Value = 'Litecoin'
$ws.Range('C19').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D19').Value = '69.39'
$ws.Range('E19').Value = '  -0.25%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').Value = '5.91'
$ws.Range('E20').Value = '  -1.44%  '
$ws.Range('D21').Value = '0.0₃0820'
$ws.Range('E21').Value = '  -0.08%  '
$ws.Range('D22').Value = '223.53'
$ws.Range('E22').Value = '  -0.68%  '
$ws.Range('E23').Value = '  +0.39%  '
$ws.Range('D24').Value = '2.39'
$ws.Range('E24').Value = '  +0.40%  '
$ws.Range('D25').Value = '2.27'
$ws.Range('E25').Value = '  +2.61%  '
$ws.Range('D26').Value = '168.28'
$ws.Range('E26').Value = '  +2.19%  '
$ws.Range('D27').Value = '9.36'
$ws.Range('E27').Value = '  +1.92%  '
$ws.Range('D28').Value = '0.129'
$ws.Range('E28').Value = '  -0.85%  '
$ws.Range('D29').Value = '18.76'
$ws.Range('E29').Value = '  -0.44%  '
$ws.Range('D30').Value = '1.27'
$ws.Range('E30').Value = '  -0.57%  '
$ws.Range('D31').Value = '0.119'
$ws.Range('E31').Value = '  -0.62%  '
$ws.Range('E32').Value = '  +8.73%  '
$ws.Range('D33').Value = '4.37'
$ws.Range('E33').Value = '  -1.13%  '
$ws.Range('D34').Value = '0.0606'
$ws.Range('E34').Value = '  +0.95%  '
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('D36').Value = '6.31'
$ws.Range('E36').Value = '  -0.37%  '
$ws.Range('D37').Value = '2.32'
$ws.Range('E37').Value = '  +3.52%  '
$ws.Range('E38').Value = '  +5.96%  '
$ws.Range('E39').Value = '  -0.11%  '
$ws.Range('D40').Value = '18.02'
$ws.Range('E40').Value = '  +9.28%  '
$ws.Range('D41').Value = '1.525.40'
$ws.Range('E41').Value = '  -0.55%  '
$ws.Range('D42').Value = '97.08'
$ws.Range('E42').Value = '  +0.76%  '
$ws.Range('D43').Value = '0.0214'
$ws.Range('E43').Value = '  -0.62%  '
$ws.Range('B44').Value = 'FTXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D44').Value = '4.35'
$ws.Range('E44').Value = '  +9.55%  '
$ws.Range('B45').Value = 'HuobiToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D45').Value = '2.84'
$ws.Range('E45').Value = '  +1.14%  '
$ws.Range('D46').Value = '0.0905'
$ws.Range('E46').Value = '  -1.34%  '
$ws.Range('D47').Value = '1.10'
$ws.Range('E47').Value = '  +0.30%  '
$ws.Range('D48').Value = '1.00'
$ws.Range('E48').Value = '  +0.49%  '
$ws.Range('E49').Value = '  -0.31%  '
$ws.Range('D50').Value = '7.01'
$ws.Range('E50').Value = '  -1.20%  '
$ws.Range('D51').Value = '2.226.09'
$ws.Range('E51').Value = '  +0.61%  '

foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "General"
}
